{"js": "// Fix the typo \"EXAMENSE SOLICITADOS:\" -> \"EXAMENES SOLICITADOS:\".\n// The original run's text is split into three runs (\"EXAMEN\" / \"ES\" /\n// \" SOLICITADOS:\") that keep the same run properties (<w:noProof/>),\n// matching how Word naturally fragments a run when only part of its\n// text is touched.\n\nconst body = context.document.body;\n\n// Locate the run containing the misspelled label.\nconst results = body.search(\"EXAMENSE SOLICITADOS:\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n\n  // Narrow down to the two-letter span that needs to be reordered\n  // (\"SE\" -> \"ES\"); it only occurs once inside the matched text.\n  const sub = hit.search(\"SE\", { matchCase: true });\n  sub.load(\"items\");\n  await context.sync();\n\n  if (sub.items.length > 0) {\n    const target = sub.items[0];\n\n    // Correct the text in place.\n    target.insertText(\"ES\", \"Replace\");\n\n    // Toggling a character property forces Word to materialize this\n    // sub-range as its own run (splitting the original run into three),\n    // then reverting the property keeps the formatting identical to the\n    // surrounding text while leaving the run split in place.\n    target.font.set({ bold: true });\n    target.font.set({ bold: false });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix the typo \"EXAMENSE SOLICITADOS:\" -> \"EXAMENES SOLICITADOS:\".\n# The original run's text is split into three runs (\"EXAMEN\" / \"ES\" /\n# \" SOLICITADOS:\") that keep the same run properties (<w:noProof/>),\n# matching how Word naturally fragments a run when only part of its\n# text is touched.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph/run containing the misspelled label.\n$hit = $d.Content\n$hit.Find.Text = \"EXAMENSE SOLICITADOS:\"\n\nif ($hit.Find.Execute()) {\n    # Narrow down to the two-letter span that needs to be reordered\n    # (\"SE\" -> \"ES\"); scoped to the match so it can't find it elsewhere.\n    $sub = $hit.Duplicate\n    $sub.Find.Text = \"SE\"\n\n    if ($sub.Find.Execute()) {\n        # Correct the text in place.\n        $sub.Text = \"ES\"\n\n        # Toggling a character property forces Word to materialize this\n        # sub-range as its own run (splitting the original run into three),\n        # then reverting the property keeps the formatting identical to the\n        # surrounding text while leaving the run split in place.\n        $sub.Font.Bold = 1\n        $sub.Font.Bold = 0\n    }\n}\n"}
